$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 9853.63442723393
$ws.Range("C2").Value = 8985.72695820076
$ws.Range("E2").Value = 5123.32585128221
$ws.Range("F2").Value = -8.26446627154284
$ws.Range("B3").Value = 10328.1275783307
$ws.Range("C3").Value = 9776.47127469568
$ws.Range("E3").Value = 5545.29217689646
$ws.Range("F3").Value = 263.265143816339
$ws.Range("B4").Value = 10440.1848935311
$ws.Range("C4").Value = 9769.31085332607
$ws.Range("E4").Value = 6034.38003079527
$ws.Range("F4").Value = 283.345453505056
$ws.Range("B5").Value = 10448.497791708
$ws.Range("C5").Value = 9411.42474045777
$ws.Range("E5").Value = 6060.24320125144
$ws.Range("F5").Value = 269.511164237883
$ws.Range("B6").Value = 4423.76645036624
$ws.Range("C6").Value = 7147.28802982712
$ws.Range("E6").Value = 5992.31773317116
$ws.Range("F6").Value = 172.341906791595
$ws.Range("B7").Value = 4536.89570607893
$ws.Range("C7").Value = 7366.71785145019
$ws.Range("E7").Value = 6312.26773984401
$ws.Range("F7").Value = 194.816066303925
$ws.Range("C9").Value = 10207.4075655513
$ws.Range("F9").Value = 317.217851167791
$ws.Range("C10").Value = 10476.1323134924
$ws.Range("F10").Value = 328.414715665339
$ws.Range("C11").Value = 10751.558110052
$ws.Range("F11").Value = 339.890790521986
$ws.Range("C12").Value = 10276.5604829325
$ws.Range("F12").Value = 320.099222725344
$ws.Range("C13").Value = 7511.29884228622
$ws.Range("F13").Value = 189.204138817247
$ws.Range("C14").Value = 7229.32092561115
$ws.Range("F14").Value = 177.10136804071
$ws.Range("C15").Value = 11099.9535410377
$ws.Range("F15").Value = 408.500307435758
